# "update security filter, rename orderdetails entity, rename permission"
#
# The workbook holds a single-column list of permission codes in Sheet1!A.
# Several permission names are renamed (e.g. CREATE_ORDERDETAIL ->
# CREATE_ORDER_DETAILS, VIEW_ACCOUNTS -> VIEW_ACCOUNT, ...), and the
# (unchanged) alphabetical order of the list is preserved.
#
# The final, authoritative list of values for A1:A41 after the edit:

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "name",
    "CREATE_ACCOUNT",
    "CREATE_ACCOUNT_ROLE_PERMISSION",
    "CREATE_CUSTOMER",
    "CREATE_ORDER",
    "CREATE_ORDER_DETAILS",
    "CREATE_ORDER_STATUS",
    "CREATE_PERMISSION",
    "CREATE_PRODUCT",
    "CREATE_PRODUCT_TYPE",
    "CREATE_ROLE",
    "DELETE_ACCOUNT",
    "DELETE_ACCOUNT_ROLE_PERMISSION",
    "DELETE_CUSTOMER",
    "DELETE_ORDER",
    "DELETE_ORDER_DETAILS",
    "DELETE_ORDER_STATUS",
    "DELETE_PERMISSION",
    "DELETE_PRODUCT",
    "DELETE_PRODUCT_TYPE",
    "DELETE_ROLE",
    "UPDATE_ACCOUNT",
    "UPDATE_ACCOUNT_ROLE_PERMISSION",
    "UPDATE_CUSTOMER",
    "UPDATE_ORDER",
    "UPDATE_ORDER_DETAILS",
    "UPDATE_ORDER_STATUS",
    "UPDATE_PERMISSION",
    "UPDATE_PRODUCT",
    "UPDATE_PRODUCT_TYPE",
    "UPDATE_ROLE",
    "VIEW_ACCOUNT",
    "VIEW_ACCOUNT_ROLE_PERMISSION",
    "VIEW_CUSTOMER",
    "VIEW_ORDER",
    "VIEW_ORDER_DETAILS",
    "VIEW_ORDER_STATUS",
    "VIEW_PERMISSION",
    "VIEW_PRODUCT",
    "VIEW_PRODUCT_TYPE",
    "VIEW_ROLE"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

# Update the selection / scroll position recorded in the sheet view to match
# the commit (cursor left at E38, near the bottom of the renamed list).
$ws.Range("E38").Select()
